$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 7 (pushes the old "Desc" row down to row 9)
$ws.Range("A7:A8").EntireRow.Insert()

# Row 7: new "Icon" field
$ws.Range("A7").Value = "Icon"
$ws.Range("B7").Value = "string"
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = $false
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = "Friend"
$ws.Range("J7").Value = "图标"

# Row 8: new "ShowName" field
$ws.Range("A8").Value = "ShowName"
$ws.Range("B8").Value = "string"
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = $false
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = "Friend"
$ws.Range("J8").Value = "名字"

# Move selection cursor to C13
$ws.Range("C13").Select()
